$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add three new rows (11-13) for the new "NewShape3" test group ---
# Set column A (Test Name) values first, then column B (Test Group) values,
# so new shared-string entries are interned in the same order as the source
# edit (NSP3.1, NSP3.2, NSP3.3, then NewShape3).
$ws.Range("A11").Value = "NSP3.1"
$ws.Range("A12").Value = "NSP3.2"
$ws.Range("A13").Value = "NSP3.3"
$ws.Range("B11").Value = "NewShape3"
$ws.Range("B12").Value = "NewShape3"
$ws.Range("B13").Value = "NewShape3"

# Match the number formats used by the rest of the table for these columns
# (D = Percent style, E/F = 0.00 numeric style) before writing the formulas.
$ws.Range("D11:D13").NumberFormat = "0.00%"
$ws.Range("E11:F13").NumberFormat = "0.00"

$ws.Range("D11").Formula = "=1/3"
$ws.Range("D12").Formula = "=1/3"
$ws.Range("D13").Formula = "=1/3"

$ws.Range("E11").Formula = "=C11/D11"
$ws.Range("E12").Formula = "=C12/D12"
$ws.Range("E13").Formula = "=C13/D13"

$ws.Range("F11").Formula = "=AVERAGEIF(B:B,B11,E:E)"
$ws.Range("F12").Formula = "=AVERAGEIF(B:B,B12,E:E)"
$ws.Range("F13").Formula = "=AVERAGEIF(B:B,B13,E:E)"

# --- Rebuild the three conditional-formatting rules so Excel mints fresh
# dxf style records (as happens whenever these rules get re-touched),
# keeping the same ranges, colors, order and priorities as before. ---
$ws.Range("E1:F1048576").FormatConditions.Delete()
$ws.Range("F1:F1048576").FormatConditions.Delete()
$ws.Range("C1:C1048576").FormatConditions.Delete()

$fcChange = $ws.Range("E1:F1048576").FormatConditions.AddTop10()
$fcChange.TopBottom = 5
$fcChange.Rank = 10
$fcChange.Percent = $true
$fcChange.Font.Color = 24832
$fcChange.Interior.Color = 13561798
$fcChange.Priority = 3

$fcGroupChange = $ws.Range("F1:F1048576").FormatConditions.AddTop10()
$fcGroupChange.TopBottom = 5
$fcGroupChange.Rank = 10
$fcGroupChange.Percent = $true
$fcGroupChange.Font.Color = 24832
$fcGroupChange.Interior.Color = 13561798
$fcGroupChange.Priority = 2

$fcAccuracy = $ws.Range("C1:C1048576").FormatConditions.Add(1, 5, "95")
$fcAccuracy.Font.Color = 24832
$fcAccuracy.Interior.Color = 13561798
$fcAccuracy.Priority = 1

# --- Move the active selection to match where the user left off editing ---
$ws.Range("D17").Select() | Out-Null

# --- Restore the saved window position of the workbook, best-effort ---
$w = $excel.ActiveWindow
$w.Left = -20
$w.Top = 460
